# Weekly data refresh: insert a new "Apio" price entry at the top of the
# time series (row 169), pushing the existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 169 - everything below (169..257) shifts down to 170..258
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Cells.Item(169, 1).Value = 7
$ws.Cells.Item(169, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(169, 3).Value = "Ñuble"
$ws.Cells.Item(169, 4).Value = 44879
$ws.Cells.Item(169, 5).Value = 16
$ws.Cells.Item(169, 6).Value = 100112017
$ws.Cells.Item(169, 7).Value = "Apio"
$ws.Cells.Item(169, 8).Value = "Americana (o)"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 120
$ws.Cells.Item(169, 11).Value = 7500
$ws.Cells.Item(169, 12).Value = 8000
$ws.Cells.Item(169, 13).Value = 7750
$ws.Cells.Item(169, 14).Value = "$/docena de matas"
$ws.Cells.Item(169, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(169, 16).Value = 1292
$ws.Cells.Item(169, 17).Value = 6
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# Match the date column's existing date/time number format (style index 2
# in the original workbook) so the new row's "Fecha" cell renders the same
# way as the rest of column D.
$ws.Cells.Item(169, 4).NumberFormat = $ws.Cells.Item(170, 4).NumberFormat
